# Auto-generated cell updates derived from the OOXML diff.
# Updates Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# with refreshed market-price figures pulled by the scheduled price-scrape runner.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7666.6665   # was 8000
$ws.Range("I62").Value = 6500   # was 7000
$ws.Range("K62").Value = 6500   # was 7000
$ws.Range("M62").Value = -5876   # was -6376
$ws.Range("H65").Value = 7666.6665   # was 8000
$ws.Range("I65").Value = 6500   # was 7000
$ws.Range("K65").Value = 32500   # was 35000
$ws.Range("M65").Value = -29380   # was -31880
$ws.Range("H121").Value = 1070.1666   # was 951.6667
$ws.Range("J121").Value = 1031.091   # was 912.5
$ws.Range("L121").Value = 3093.273   # was 2737.5
$ws.Range("N121").Value = -6587.272999999999   # was -6231.5
$ws.Range("H137").Value = 1462.88   # was 1510.0869
$ws.Range("I137").Value = 1016.2   # was 1031
$ws.Range("K137").Value = 3048.6   # was 3093
$ws.Range("M137").Value = -498.6000000000004   # was -543
$ws.Range("H138").Value = 1577.54   # was 1446.7959
$ws.Range("I138").Value = 1297.3125   # was 878.3461
$ws.Range("J138").Value = 1709.4117   # was 1652.0695
$ws.Range("K138").Value = 3891.9375   # was 2635.0383
$ws.Range("L138").Value = 5128.2351   # was 4956.208500000001
$ws.Range("M138").Value = 1248.0625   # was 2504.9617
$ws.Range("N138").Value = -15408.2351   # was -15236.2085
$ws.Range("H141").Value = 612.8570999999999   # was 615.3570999999999
$ws.Range("I141").Value = 612.8570999999999   # was 615.3570999999999
$ws.Range("K141").Value = 1838.5713   # was 1846.0713
$ws.Range("M141").Value = 3341.4287   # was 3333.9287

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 4528.4614   # was 4362.143
$ws.Range("I28").Value = 4528.4614   # was 4362.143
$ws.Range("K28").Value = 4528.4614   # was 4362.143
$ws.Range("M28").Value = -4336.4614   # was -4170.143
$ws.Range("H32").Value = 3469.4304   # was 3553.026
$ws.Range("I32").Value = 3160.5278   # was 3243.6572
$ws.Range("K32").Value = 3160.5278   # was 3243.6572
$ws.Range("M32").Value = -2873.5278   # was -2956.6572
$ws.Range("H45").Value = 1127.0358   # was 1090.6
$ws.Range("I45").Value = 1033.8334   # was 989
$ws.Range("J45").Value = 1294.8   # was 1327.6666
$ws.Range("K45").Value = 1033.8334   # was 989
$ws.Range("L45").Value = 1294.8   # was 1327.6666
$ws.Range("M45").Value = -656.8334   # was -612
$ws.Range("N45").Value = -2048.8   # was -2081.6666
$ws.Range("H74").Value = 1612.862   # was 1550.4375
$ws.Range("I74").Value = 876.7222   # was 892.8333
$ws.Range("J74").Value = 2817.4546   # was 2395.9285
$ws.Range("K74").Value = 876.7222   # was 892.8333
$ws.Range("L74").Value = 2817.4546   # was 2395.9285
$ws.Range("M74").Value = -2.722200000000043   # was -18.83330000000001
$ws.Range("N74").Value = -4565.4546   # was -4143.9285
$ws.Range("H77").Value = 1612.862   # was 1550.4375
$ws.Range("I77").Value = 876.7222   # was 892.8333
$ws.Range("J77").Value = 2817.4546   # was 2395.9285
$ws.Range("K77").Value = 4383.611   # was 4464.1665
$ws.Range("L77").Value = 14087.273   # was 11979.6425
$ws.Range("M77").Value = -15.61099999999988   # was -96.16650000000027
$ws.Range("N77").Value = -22823.273   # was -20715.6425
$ws.Range("H99").Value = 4528.4614   # was 4362.143
$ws.Range("I99").Value = 4528.4614   # was 4362.143
$ws.Range("K99").Value = 4528.4614   # was 4362.143
$ws.Range("M99").Value = -1533.4614   # was -1367.143
$ws.Range("H122").Value = 976.5454999999999   # was 1028
$ws.Range("I122").Value = 1002.8   # was 1062.8889
$ws.Range("K122").Value = 3008.4   # was 3188.6667
$ws.Range("M122").Value = -558.3999999999996   # was -738.6666999999998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1500   # was 0
$ws.Range("I12").Value = 1500   # was 0
$ws.Range("K12").Value = 1500   # was 0
$ws.Range("M12").Value = -1332   # was None
$ws.Range("H105").Value = 91809730   # was 100990530
$ws.Range("I105").Value = 100990560   # was 112211544
$ws.Range("K105").Value = 100990560   # was 112211544
$ws.Range("M105").Value = -100988813   # was -112209797
$ws.Range("H134").Value = 5877.931   # was 5894.1377
$ws.Range("I134").Value = 1667.35   # was 1690.85
$ws.Range("K134").Value = 5002.049999999999   # was 5072.549999999999
$ws.Range("M134").Value = -2467.049999999999   # was -2537.549999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1847   # was 1566.3334
$ws.Range("I31").Value = 1629.3334   # was 1499.5
$ws.Range("J31").Value = 2500   # was 1700
$ws.Range("K31").Value = 1629.3334   # was 1499.5
$ws.Range("L31").Value = 2500   # was 1700
$ws.Range("M31").Value = -1334.3334   # was -1204.5
$ws.Range("N31").Value = -3090   # was -2290
$ws.Range("H34").Value = 1847   # was 1566.3334
$ws.Range("I34").Value = 1629.3334   # was 1499.5
$ws.Range("J34").Value = 2500   # was 1700
$ws.Range("K34").Value = 1629.3334   # was 1499.5
$ws.Range("L34").Value = 2500   # was 1700
$ws.Range("M34").Value = -1427.3334   # was -1297.5
$ws.Range("N34").Value = -2904   # was -2104
$ws.Range("H86").Value = 3522766   # was 3346767.2
$ws.Range("J86").Value = 23846.2   # was 21932
$ws.Range("L86").Value = 23846.2   # was 21932
$ws.Range("N86").Value = -26092.2   # was -24178
$ws.Range("H89").Value = 3522766   # was 3346767.2
$ws.Range("J89").Value = 23846.2   # was 21932
$ws.Range("L89").Value = 119231   # was 109660
$ws.Range("N89").Value = -130463   # was -120892
$ws.Range("H107").Value = 614.2917   # was 650.0476
$ws.Range("I107").Value = 563.6667   # was 547.1539
$ws.Range("J107").Value = 698.6667   # was 817.25
$ws.Range("K107").Value = 563.6667   # was 547.1539
$ws.Range("L107").Value = 698.6667   # was 817.25
$ws.Range("M107").Value = 1356.3333   # was 1372.8461
$ws.Range("N107").Value = -4538.6667   # was -4657.25
$ws.Range("H114").Value = 23921.54   # was 24796
$ws.Range("I114").Value = 0   # was 21000
$ws.Range("J114").Value = 23921.54   # was 25745
$ws.Range("K114").Value = 0   # was 21000
$ws.Range("L114").Value = 23921.54   # was 25745
$ws.Range("M114").Value = $null   # was -16661
$ws.Range("N114").Value = -32599.54   # was -34423
$ws.Range("H132").Value = 1632.7587   # was 1422.7428
$ws.Range("I132").Value = 925.0909   # was 814.2143
$ws.Range("K132").Value = 2775.2727   # was 2442.6429
$ws.Range("M132").Value = -245.2727   # was 87.35710000000017
$ws.Range("H134").Value = 881   # was 948.1429000000001
$ws.Range("I134").Value = 720.4706   # was 786.7857
$ws.Range("K134").Value = 2161.4118   # was 2360.3571
$ws.Range("M134").Value = 373.5882000000001   # was 174.6428999999998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 291.94736   # was 288.22223
$ws.Range("I92").Value = 436.75   # was 387.6
$ws.Range("J92").Value = 253.33333   # was 250
$ws.Range("K92").Value = 1310.25   # was 1162.8
$ws.Range("L92").Value = 759.99999   # was 750
$ws.Range("M92").Value = -62.25   # was 85.19999999999982
$ws.Range("N92").Value = -3255.99999   # was -3246
$ws.Range("H102").Value = 2957.8   # was 2941.5
$ws.Range("J102").Value = 2957.8   # was 2941.5
$ws.Range("L102").Value = 8873.400000000001   # was 8824.5
$ws.Range("N102").Value = -13741.4   # was -13692.5
$ws.Range("H109").Value = 79361.53999999999   # was 68817.336
$ws.Range("I109").Value = 167283.33   # was 125532.5
$ws.Range("K109").Value = 501849.99   # was 376597.5
$ws.Range("M109").Value = -500809.99   # was -375557.5
$ws.Range("H113").Value = 683.78125   # was 685.8484999999999
$ws.Range("J113").Value = 689.43335   # was 691.4516
$ws.Range("L113").Value = 2068.30005   # was 2074.3548
$ws.Range("N113").Value = -6408.30005   # was -6414.3548
$ws.Range("H131").Value = 10418694   # was 12822812
$ws.Range("J131").Value = 2143.6223   # was 2467.0278
$ws.Range("L131").Value = 6430.8669   # was 7401.0834
$ws.Range("N131").Value = -16510.8669   # was -17481.0834

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 5011   # was 0
$ws.Range("J27").Value = 5011   # was 0
$ws.Range("L27").Value = 5011   # was 0
$ws.Range("N27").Value = -5343   # was None
$ws.Range("H28").Value = 0   # was 5000
$ws.Range("J28").Value = 0   # was 5000
$ws.Range("L28").Value = 0   # was 5000
$ws.Range("N28").Value = $null   # was -5384
$ws.Range("H102").Value = 3501.8108   # was 3730.647
$ws.Range("I102").Value = 2329.7693   # was 2515.1738
$ws.Range("K102").Value = 2329.7693   # was 2515.1738
$ws.Range("M102").Value = -707.7692999999999   # was -893.1738
$ws.Range("H122").Value = 2818.7   # was 3042.7144
$ws.Range("I122").Value = 2423.5   # was 2716.6667
$ws.Range("J122").Value = 4399.5   # was 4999
$ws.Range("K122").Value = 7270.5   # was 8150.000100000001
$ws.Range("L122").Value = 13198.5   # was 14997
$ws.Range("M122").Value = -4820.5   # was -5700.000100000001
$ws.Range("N122").Value = -18098.5   # was -19897

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 40000   # was 0
$ws.Range("J125").Value = 40000   # was 0
$ws.Range("L125").Value = 40000   # was 0
$ws.Range("N125").Value = -49840   # was None
$ws.Range("H132").Value = 22099.49   # was 19833.164
$ws.Range("I132").Value = 1084.2222   # was 1067.1333
$ws.Range("J132").Value = 47890.953   # was 42352.4
$ws.Range("K132").Value = 3252.6666   # was 3201.3999
$ws.Range("L132").Value = 143672.859   # was 127057.2
$ws.Range("M132").Value = -722.6665999999996   # was -671.3998999999999
$ws.Range("N132").Value = -148732.859   # was -132117.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1660.375   # was 1870.381
$ws.Range("I132").Value = 1452.5555   # was 1705
$ws.Range("K132").Value = 4357.666499999999   # was 5115
$ws.Range("M132").Value = -1827.666499999999   # was -2585
$ws.Range("H136").Value = 569   # was 662.44446
$ws.Range("I136").Value = 270   # was 310.5
$ws.Range("K136").Value = 810   # was 931.5
$ws.Range("M136").Value = 1740   # was 1618.5
